$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Quantity Sold (C), Expenses (E) and Financial Result (F) values
# for each weapon row, then the Totals row (12).

$ws.Cells.Item(2, 3).Value = 5
$ws.Cells.Item(2, 5).Value = 10
$ws.Cells.Item(2, 6).Value = 290

$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 4999

$ws.Cells.Item(4, 3).Value = 1843
$ws.Cells.Item(4, 5).Value = 1843
$ws.Cells.Item(4, 6).Value = 6957

$ws.Cells.Item(5, 3).Value = 1298
$ws.Cells.Item(5, 5).Value = 1298
$ws.Cells.Item(5, 6).Value = 21452

$ws.Cells.Item(6, 3).Value = 1140
$ws.Cells.Item(6, 5).Value = 2280
$ws.Cells.Item(6, 6).Value = 20820

$ws.Cells.Item(7, 3).Value = 42
$ws.Cells.Item(7, 5).Value = 84
$ws.Cells.Item(7, 6).Value = 25116

$ws.Cells.Item(8, 3).Value = 22
$ws.Cells.Item(8, 5).Value = 44
$ws.Cells.Item(8, 6).Value = 32956

$ws.Cells.Item(9, 3).Value = 2017
$ws.Cells.Item(9, 5).Value = 2017
$ws.Cells.Item(9, 6).Value = 61983

$ws.Cells.Item(10, 3).Value = 2108
$ws.Cells.Item(10, 5).Value = 4216
$ws.Cells.Item(10, 6).Value = 70584

$ws.Cells.Item(11, 3).Value = 100
$ws.Cells.Item(11, 5).Value = 100
$ws.Cells.Item(11, 6).Value = 1999900

$ws.Cells.Item(12, 3).Value = 8576
$ws.Cells.Item(12, 5).Value = 11893
$ws.Cells.Item(12, 6).Value = 2245057

$wb.Save()
